# AdornmentsEULA.docx edit:
#  - The placeholder copyright line "Copyright (c) <year> <copyright holders>"
#    is replaced with the concrete "Copyright (c) 2015 Microsoft Corp."
#  - The hidden "_GoBack" bookmark (which Word silently re-drops at the most
#    recent edit point on every save) moves from the end of the "THE SOFTWARE
#    IS PROVIDED..." paragraph to the end of the now-edited copyright
#    paragraph.

$d = $word.ActiveDocument

# 1) Replace the MIT license placeholder copyright line with the real one.
$d.Content.Find.Execute(
    "Copyright (c) <year> <copyright holders>", # FindText
    $true,                                       # MatchCase
    $false,                                      # MatchWholeWord
    $false,                                      # MatchWildcards
    $false,                                      # MatchSoundsLike
    $false,                                      # MatchAllWordForms
    $true,                                       # Forward
    1,                                            # Wrap (wdFindContinue)
    $false,                                      # Format
    "Copyright (c) 2015 Microsoft Corp.",        # ReplaceWith
    2                                             # Replace (wdReplaceAll)
) | Out-Null

# 2) Move the "_GoBack" bookmark onto the edited paragraph. Adding a bookmark
#    with a name that already exists elsewhere replaces/moves it, exactly as
#    Word does when you type - the bookmark always tracks the last edit.
$copyrightPara = $d.Paragraphs(2).Range
$copyrightPara.End = $copyrightPara.End - 1
$d.Bookmarks.Add("_GoBack", $copyrightPara) | Out-Null
